$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 (B6:F6) previously held date serial numbers formatted to look like
# "1945年8月15日" etc. They are replaced with literal text strings that
# include the quote characters themselves.
$ws.Range("B6").Value = '"1945年8月15日"'
$ws.Range("C6").Value = '"1949年10月1日"'
$ws.Range("D6").Value = '"1950年10月1日"'
$ws.Range("E6").Value = '"1949年9月21日"'
$ws.Range("F6").Value = '"1950年1月1日"'

# Update the view: select K19 instead of Z8, and drop the E1 scroll anchor
# (scroll back to the top-left of the sheet).
$ws.Activate()
$ws.Range("K19").Select()
